# "Bento object repository revisited"
# The FilesTab Neo4j/Cypher query stored in cell B4 of the "startup" sheet
# drops the `File Type` and `Breed` projected columns, and the sheet's
# saved selection moves from D3 to B4 (scrolled so row 4 is in view).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 'MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
 MATCH (samp:sample)-->(c) 
 WHERE samp.specific_sample_pathology IN ["Lymphoma"]  
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '''') AS `File Name`, 
        coalesce(labels(parent)[0], '''') AS `Association`,
        coalesce(f.file_description, '''') AS `Description`,
        coalesce(f.file_format, '''') AS `Format`,
        coalesce(f.file_size, '''') AS `Size`,
        coalesce(c.case_id, '''') AS `Case ID`, 
        coalesce(diag.disease_term,'''') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'''') AS `Study Code`'

# Reflect the new active cell / scroll position recorded in the workbook.
$ws.Range("B4").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
